# Insert a new data row at row 678 (2026/01/22, 木, 19, 20), pushing the
# existing rows 678-719 down to 679-720.
#
# Row 677 already holds ("2026/01/22", "木", 16, 20) - i.e. the same date
# and weekday text as the new row, just a different hour/rank. Copying it
# down (instead of retyping the date string) keeps column A as the plain
# text it already is, rather than letting Excel's date auto-detection turn
# a freshly-typed "2026/01/22" into a date serial value/format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("677").Copy()
$ws.Rows("678").Insert()

$ws.Range("C678").Value = 19
